# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-market-day records (rows 2-13)
# across the existing rows. Columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen) and
# S (Precio $/Kg) are the ones that carry the per-row data that moves;
# all the other columns stay constant for every row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) whose values travel together with each record.
$cols = @("D", "M", "N", "O", "P", "R", "S")

# Snapshot the current ("before") values for rows 2-13 so we can re-write
# them in their new order without clobbering data we still need to read.
$snapshot = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Mapping: new row number -> old row number whose values now land there.
$mapping = @{
    2  = 6
    3  = 12
    4  = 10
    5  = 8
    6  = 9
    7  = 2
    8  = 13
    9  = 11
    10 = 7
    11 = 4
    12 = 3
    13 = 5
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $source = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $source[$c]
    }
}
